$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sending-cluster = "MuSCs" rows (previously rows 8-10),
# which are dropped entirely now that "ECs" is no longer a sending cluster
# and the remaining combinations shift up.
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Rewrite rows 2-7 with the updated TPM-derived values.
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.583520999999999
$ws.Range("H2").Value = 25.750563
$ws.Range("I2").Value = 0.8910607110509009
$ws.Range("J2").Value = 0.8910607110509009
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 69.890255385896
$ws.Range("R2").Value = 629.0122984730641
$ws.Range("S2").Value = 0.1551616431548081
$ws.Range("T2").Value = 0.1551616431548081

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.583520999999999
$ws.Range("H3").Value = 25.750563
$ws.Range("I3").Value = 0.8910607110509009
$ws.Range("J3").Value = 0.8910607110509009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 208.9258567323976
$ws.Range("R3").Value = 1880.332710591579
$ws.Range("S3").Value = 0.4638311743051196
$ws.Range("T3").Value = 0.4638311743051196

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.583520999999999
$ws.Range("H4").Value = 25.750563
$ws.Range("I4").Value = 0.8910607110509009
$ws.Range("J4").Value = 0.8910607110509009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 122.548937860914
$ws.Range("R4").Value = 1102.940440748226
$ws.Range("S4").Value = 0.2720678935909732
$ws.Range("T4").Value = 0.2720678935909732

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.049404
$ws.Range("H5").Value = 3.148212
$ws.Range("I5").Value = 0.1089392889490991
$ws.Range("J5").Value = 0.1089392889490991
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 8.544641943904001
$ws.Range("R5").Value = 76.90177749513602
$ws.Range("S5").Value = 0.01896975017282864
$ws.Range("T5").Value = 0.01896975017282864

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.049404
$ws.Range("H6").Value = 3.148212
$ws.Range("I6").Value = 0.1089392889490991
$ws.Range("J6").Value = 0.1089392889490991
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("Q6").Value = 25.54285470477733
$ws.Range("R6").Value = 229.885692342996
$ws.Range("S6").Value = 0.05670706574149347
$ws.Range("T6").Value = 0.05670706574149347

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.049404
$ws.Range("H7").Value = 3.148212
$ws.Range("I7").Value = 0.1089392889490991
$ws.Range("J7").Value = 0.1089392889490991
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("Q7").Value = 14.982586468536
$ws.Range("R7").Value = 134.843278216824
$ws.Range("S7").Value = 0.03326247303477695
$ws.Range("T7").Value = 0.03326247303477695

